$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 91, pushing existing rows 91..164 down to 93..166.
$ws.Rows("91:92").Insert()

# New row 91: Femacal de La Calera / Coquimbo / Alcachofa / Argentina(o) / Primera
$ws.Cells.Item(91,1).Value  = 3
$ws.Cells.Item(91,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(91,3).Value  = "Coquimbo"
$ws.Cells.Item(91,4).Value  = [DateTime]"2021-09-08"
$ws.Cells.Item(91,5).Value  = 5
$ws.Cells.Item(91,6).Value  = 100112013
$ws.Cells.Item(91,7).Value  = "Alcachofa"
$ws.Cells.Item(91,8).Value  = "Argentina(o)"
$ws.Cells.Item(91,9).Value  = "Primera"
$ws.Cells.Item(91,10).Value = 68
$ws.Cells.Item(91,11).Value = 10000
$ws.Cells.Item(91,12).Value = 11000
$ws.Cells.Item(91,13).Value = 10559
$ws.Cells.Item(91,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(91,15).Value = "Provincia de Limarí"
$ws.Cells.Item(91,16).Value = 211
$ws.Cells.Item(91,17).Value = 50
$ws.Cells.Item(91,18).Value = "Hortaliza"

# New row 92: Femacal de La Calera / Coquimbo / Alcachofa / Española / Extra
$ws.Cells.Item(92,1).Value  = 3
$ws.Cells.Item(92,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(92,3).Value  = "Coquimbo"
$ws.Cells.Item(92,4).Value  = [DateTime]"2021-09-08"
$ws.Cells.Item(92,5).Value  = 5
$ws.Cells.Item(92,6).Value  = 100112013
$ws.Cells.Item(92,7).Value  = "Alcachofa"
$ws.Cells.Item(92,8).Value  = "Española"
$ws.Cells.Item(92,9).Value  = "Extra"
$ws.Cells.Item(92,10).Value = 73
$ws.Cells.Item(92,11).Value = 12000
$ws.Cells.Item(92,12).Value = 12500
$ws.Cells.Item(92,13).Value = 12260
$ws.Cells.Item(92,14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(92,15).Value = "Provincia de Limarí"
$ws.Cells.Item(92,16).Value = 409
$ws.Cells.Item(92,17).Value = 30
$ws.Cells.Item(92,18).Value = "Hortaliza"
